$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.243977666666667
$ws.Range("H2").Value = 9.731933000000001
$ws.Range("I2").Value = 0.0124341611854976
$ws.Range("J2").Value = 0.01249584677475898
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.084959
$ws.Range("N2").Value = 63.25487699999999
$ws.Range("O2").Value = 0.03381320693734752
$ws.Range("P2").Value = 0.03509122472428063
$ws.Range("Q2").Value = 68.39913609858235
$ws.Range("R2").Value = 615.5922248872411
$ws.Range("S2").Value = 0.0004204388652575646
$ws.Range("T2").Value = 0.0004384945672932445

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.243977666666667
$ws.Range("H3").Value = 9.731933000000001
$ws.Range("I3").Value = 0.0124341611854976
$ws.Range("J3").Value = 0.01249584677475898
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 81.06331633333333
$ws.Range("N3").Value = 243.189949
$ws.Range("O3").Value = 0.12999838843446
$ws.Range("P3").Value = 0.1349118606466557
$ws.Range("Q3").Value = 262.9675877712686
$ws.Range("R3").Value = 2366.708289941418
$ws.Range("S3").Value = 0.001616420915649002
$ws.Range("T3").Value = 0.001685837938738245

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.243977666666667
$ws.Range("H4").Value = 9.731933000000001
$ws.Range("I4").Value = 0.0124341611854976
$ws.Range("J4").Value = 0.01249584677475898
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 136.9994176666667
$ws.Range("N4").Value = 410.998253
$ws.Range("O4").Value = 0.2197011461990087
$ws.Range("P4").Value = 0.2280050605000741
$ws.Range("Q4").Value = 444.4230512570055
$ws.Range("R4").Value = 3999.807461313049
$ws.Range("S4").Value = 0.002731799464477047
$ws.Range("T4").Value = 0.002849116299878576

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.243977666666667
$ws.Range("H5").Value = 9.731933000000001
$ws.Range("I5").Value = 0.0124341611854976
$ws.Range("J5").Value = 0.01249584677475898
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 316.292811
$ws.Range("N5").Value = 948.878433
$ws.Range("O5").Value = 0.5072276531881493
$ws.Range("P5").Value = 0.5263990368430604
$ws.Range("Q5").Value = 1026.046815011221
$ws.Range("R5").Value = 9234.421335100989
$ws.Range("S5").Value = 0.006306950397483124
$ws.Range("T5").Value = 0.006577801706771589

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.243977666666667
$ws.Range("H6").Value = 9.731933000000001
$ws.Range("I6").Value = 0.0124341611854976
$ws.Range("J6").Value = 0.01249584677475898
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 68.131198
$ws.Range("N6").Value = 136.262396
$ws.Range("O6").Value = 0.1092596052410345
$ws.Range("P6").Value = 0.07559281728592908
$ws.Range("Q6").Value = 221.0160847152447
$ws.Range("R6").Value = 1326.096508291468
$ws.Range("S6").Value = 0.001358551542630861
$ws.Range("T6").Value = 0.0009445962620773216

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 146.4311116666667
$ws.Range("H7").Value = 439.2933350000001
$ws.Range("I7").Value = 0.5612702158044854
$ws.Range("J7").Value = 0.5640546645083628
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 21.084959
$ws.Range("N7").Value = 63.25487699999999
$ws.Range("O7").Value = 0.03381320693734752
$ws.Range("P7").Value = 0.03509122472428063
$ws.Range("Q7").Value = 3087.493985816088
$ws.Range("R7").Value = 27787.4458723448
$ws.Range("S7").Value = 0.01897834595476677
$ws.Range("T7").Value = 0.01979336898904168

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 146.4311116666667
$ws.Range("H8").Value = 439.2933350000001
$ws.Range("I8").Value = 0.5612702158044854
$ws.Range("J8").Value = 0.5640546645083628
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 81.06331633333333
$ws.Range("N8").Value = 243.189949
$ws.Range("O8").Value = 0.12999838843446
$ws.Range("P8").Value = 0.1349118606466557
$ws.Range("Q8").Value = 11870.19152607666
$ws.Range("R8").Value = 106831.7237346899
$ws.Range("S8").Value = 0.07296422353084468
$ws.Range("T8").Value = 0.07609766429524836

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 146.4311116666667
$ws.Range("H9").Value = 439.2933350000001
$ws.Range("I9").Value = 0.5612702158044854
$ws.Range("J9").Value = 0.5640546645083628
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 136.9994176666667
$ws.Range("N9").Value = 410.998253
$ws.Range("O9").Value = 0.2197011461990087
$ws.Range("P9").Value = 0.2280050605000741
$ws.Range("Q9").Value = 20060.97702661597
$ws.Range("R9").Value = 180548.7932395438
$ws.Range("S9").Value = 0.1233117097396104
$ws.Range("T9").Value = 0.1286073179065783

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 146.4311116666667
$ws.Range("H10").Value = 439.2933350000001
$ws.Range("I10").Value = 0.5612702158044854
$ws.Range("J10").Value = 0.5640546645083628
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 316.292811
$ws.Range("N10").Value = 948.878433
$ws.Range("O10").Value = 0.5072276531881493
$ws.Range("P10").Value = 0.5263990368430604
$ws.Range("Q10").Value = 46315.10792690489
$ws.Range("R10").Value = 416835.9713421441
$ws.Range("S10").Value = 0.2846917743669153
$ws.Range("T10").Value = 0.2969178321240378

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 146.4311116666667
$ws.Range("H11").Value = 439.2933350000001
$ws.Range("I11").Value = 0.5612702158044854
$ws.Range("J11").Value = 0.5640546645083628
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 68.131198
$ws.Range("N11").Value = 136.262396
$ws.Range("O11").Value = 0.1092596052410345
$ws.Range("P11").Value = 0.07559281728592908
$ws.Range("Q11").Value = 9976.527062321777
$ws.Range("R11").Value = 59859.16237393067
$ws.Range("S11").Value = 0.06132416221234832
$ws.Range("T11").Value = 0.04263848119345669

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 5.967399666666666
$ws.Range("H12").Value = 17.902199
$ws.Range("I12").Value = 0.02287303333683595
$ws.Range("J12").Value = 0.0229865059320942
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 21.084959
$ws.Range("N12").Value = 63.25487699999999
$ws.Range("O12").Value = 0.03381320693734752
$ws.Range("P12").Value = 0.03509122472428063
$ws.Range("Q12").Value = 125.8223773082803
$ws.Range("R12").Value = 1132.401395774523
$ws.Range("S12").Value = 0.0007734106095032824
$ws.Range("T12").Value = 0.0008066246452891274

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 5.967399666666666
$ws.Range("H13").Value = 17.902199
$ws.Range("I13").Value = 0.02287303333683595
$ws.Range("J13").Value = 0.0229865059320942
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 81.06331633333333
$ws.Range("N13").Value = 243.189949
$ws.Range("O13").Value = 0.12999838843446
$ws.Range("P13").Value = 0.1349118606466557
$ws.Range("Q13").Value = 483.7372068664279
$ws.Range("R13").Value = 4353.634861797851
$ws.Range("S13").Value = 0.002973457472396352
$ws.Range("T13").Value = 0.003101152285064217

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 5.967399666666666
$ws.Range("H14").Value = 17.902199
$ws.Range("I14").Value = 0.02287303333683595
$ws.Range("J14").Value = 0.0229865059320942
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 136.9994176666667
$ws.Range("N14").Value = 410.998253
$ws.Range("O14").Value = 0.2197011461990087
$ws.Range("P14").Value = 0.2280050605000741
$ws.Range("Q14").Value = 817.530279317594
$ws.Range("R14").Value = 7357.772513858346
$ws.Range("S14").Value = 0.005025231641150994
$ws.Range("T14").Value = 0.005241039675732451

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 5.967399666666666
$ws.Range("H15").Value = 17.902199
$ws.Range("I15").Value = 0.02287303333683595
$ws.Range("J15").Value = 0.0229865059320942
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 316.292811
$ws.Range("N15").Value = 948.878433
$ws.Range("O15").Value = 0.5072276531881493
$ws.Range("P15").Value = 0.5263990368430604
$ws.Range("Q15").Value = 1887.445614930463
$ws.Range("R15").Value = 16987.01053437417
$ws.Range("S15").Value = 0.0116018350207376
$ws.Range("T15").Value = 0.01210007458304168

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 5.967399666666666
$ws.Range("H16").Value = 17.902199
$ws.Range("I16").Value = 0.02287303333683595
$ws.Range("J16").Value = 0.0229865059320942
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 68.131198
$ws.Range("N16").Value = 136.262396
$ws.Range("O16").Value = 0.1092596052410345
$ws.Range("P16").Value = 0.07559281728592908
$ws.Range("Q16").Value = 406.5660882348006
$ws.Range("R16").Value = 2439.396529408804
$ws.Range("S16").Value = 0.002499098593047718
$ws.Range("T16").Value = 0.001737614742966722

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 101.3861963333333
$ws.Range("H17").Value = 304.158589
$ws.Range("I17").Value = 0.388613127687944
$ws.Range("J17").Value = 0.3905410285264901
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 21.084959
$ws.Range("N17").Value = 63.25487699999999
$ws.Range("O17").Value = 0.03381320693734752
$ws.Range("P17").Value = 0.03509122472428063
$ws.Range("Q17").Value = 2137.723792854284
$ws.Range("R17").Value = 19239.51413568855
$ws.Range("S17").Value = 0.01314025610508231
$ws.Range("T17").Value = 0.01370456299607476

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 101.3861963333333
$ws.Range("H18").Value = 304.158589
$ws.Range("I18").Value = 0.388613127687944
$ws.Range("J18").Value = 0.3905410285264901
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 81.06331633333333
$ws.Range("N18").Value = 243.189949
$ws.Range("O18").Value = 0.12999838843446
$ws.Range("P18").Value = 0.1349118606466557
$ws.Range("Q18").Value = 8218.701305202439
$ws.Range("R18").Value = 73968.31174682197
$ws.Range("S18").Value = 0.05051908032390775
$ws.Range("T18").Value = 0.05268861681736742

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 101.3861963333333
$ws.Range("H19").Value = 304.158589
$ws.Range("I19").Value = 0.388613127687944
$ws.Range("J19").Value = 0.3905410285264901
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 136.9994176666667
$ws.Range("N19").Value = 410.998253
$ws.Range("O19").Value = 0.2197011461990087
$ws.Range("P19").Value = 0.2280050605000741
$ws.Range("Q19").Value = 13889.849857105
$ws.Range("R19").Value = 125008.648713945
$ws.Range("S19").Value = 0.08537874958102303
$ws.Range("T19").Value = 0.08904533083694353

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 101.3861963333333
$ws.Range("H20").Value = 304.158589
$ws.Range("I20").Value = 0.388613127687944
$ws.Range("J20").Value = 0.3905410285264901
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 316.292811
$ws.Range("N20").Value = 948.878433
$ws.Range("O20").Value = 0.5072276531881493
$ws.Range("P20").Value = 0.5263990368430604
$ws.Range("Q20").Value = 32067.72503486789
$ws.Range("R20").Value = 288609.525313811
$ws.Range("S20").Value = 0.1971153247552625
$ws.Range("T20").Value = 0.2055804212640426

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 101.3861963333333
$ws.Range("H21").Value = 304.158589
$ws.Range("I21").Value = 0.388613127687944
$ws.Range("J21").Value = 0.3905410285264901
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = 68.131198
$ws.Range("N21").Value = 136.262396
$ws.Range("O21").Value = 0.1092596052410345
$ws.Range("P21").Value = 0.07559281728592908
$ws.Range("Q21").Value = 6907.563016853207
$ws.Range("R21").Value = 41445.37810111924
$ws.Range("S21").Value = 0.0424597169226685
$ws.Range("T21").Value = 0.02952209661206178

$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 3.8636755
$ws.Range("H22").Value = 7.727351000000001
$ws.Range("I22").Value = 0.01480946198523706
$ws.Range("J22").Value = 0.00992195425829386
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 21.084959
$ws.Range("N22").Value = 63.25487699999999
$ws.Range("O22").Value = 0.03381320693734752
$ws.Range("P22").Value = 0.03509122472428063
$ws.Range("Q22").Value = 81.4654395068045
$ws.Range("R22").Value = 488.792637040827
$ws.Range("S22").Value = 0.0005007554027376021
$ws.Range("T22").Value = 0.000348173526581823

$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 3.8636755
$ws.Range("H23").Value = 7.727351000000001
$ws.Range("I23").Value = 0.01480946198523706
$ws.Range("J23").Value = 0.00992195425829386
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 81.06331633333333
$ws.Range("N23").Value = 243.189949
$ws.Range("O23").Value = 0.12999838843446
$ws.Range("P23").Value = 0.1349118606466557
$ws.Range("Q23").Value = 313.2023492658499
$ws.Range("R23").Value = 1879.214095595099
$ws.Range("S23").Value = 0.001925206191662217
$ws.Range("T23").Value = 0.001338589310237433

$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 3.8636755
$ws.Range("H24").Value = 7.727351000000001
$ws.Range("I24").Value = 0.01480946198523706
$ws.Range("J24").Value = 0.00992195425829386
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 136.9994176666667
$ws.Range("N24").Value = 410.998253
$ws.Range("O24").Value = 0.2197011461990087
$ws.Range("P24").Value = 0.2280050605000741
$ws.Range("Q24").Value = 529.3212935529672
$ws.Range("R24").Value = 3175.927761317803
$ws.Range("S24").Value = 0.003253655772747229
$ws.Range("T24").Value = 0.002262255780941259

$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 3.8636755
$ws.Range("H25").Value = 7.727351000000001
$ws.Range("I25").Value = 0.01480946198523706
$ws.Range("J25").Value = 0.00992195425829386
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 316.292811
$ws.Range("N25").Value = 948.878433
$ws.Range("O25").Value = 0.5072276531881493
$ws.Range("P25").Value = 0.5263990368430604
$ws.Range("Q25").Value = 1222.05278468683
$ws.Range("R25").Value = 7332.316708120983
$ws.Range("S25").Value = 0.007511768647750906
$ws.Range("T25").Value = 0.005222907165166791

$ws.Range("E26").Value = 2
$ws.Range("G26").Value = 3.8636755
$ws.Range("H26").Value = 7.727351000000001
$ws.Range("I26").Value = 0.01480946198523706
$ws.Range("J26").Value = 0.00992195425829386
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 68.131198
$ws.Range("N26").Value = 136.262396
$ws.Range("O26").Value = 0.1092596052410345
$ws.Range("P26").Value = 0.07559281728592908
$ws.Range("Q26").Value = 263.236840498249
$ws.Range("R26").Value = 1052.947361992996
$ws.Range("S26").Value = 0.001618075970339108
$ws.Range("T26").Value = 0.0007500284753665538

Write-Host "done"